$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.096.56"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "2.625.47"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "153.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.544"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").Value = "2.625.97"
$ws.Range("E9").Value = "  +0.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.135"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.75%  "
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000187"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.06%  "
$ws.Range("D16").Value = "3.094.85"
$ws.Range("E16").Value = "  +0.55%  "
$ws.Range("D17").Value = "68.077.66"
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("D18").Value = "2.625.11"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "367.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.39"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("E22").Value = "  -1.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("E24").Value = "  +1.00%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.87"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.754.93"
$ws.Range("E28").Value = "  +0.68%  "
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000104"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "569.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.91"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.84%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  +4.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.22%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.366"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("E42").Value = "  +2.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.95%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.83%  "
$ws.Range("D45").Value = "0.0₆0325"
$ws.Range("E45").Value = "  +10.53%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "154.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.68"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("E51").Value = "  -0.32%  "
